$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row -> source row mapping (rows whose full content must be replaced
# with the content that currently lives in the source row).
$map = @{}
$map[58] = 59
$map[59] = 58
$map[76] = 77
$map[77] = 76
$map[93] = 95
$map[94] = 93
$map[95] = 94
$map[105] = 106
$map[106] = 105
$map[113] = 114
$map[114] = 113
$map[144] = 145
$map[145] = 144
$map[150] = 151
$map[151] = 150
$map[172] = 173
$map[173] = 172
$map[174] = 175
$map[175] = 174
$map[181] = 182
$map[182] = 181
$map[186] = 187
$map[187] = 186
$map[190] = 191
$map[191] = 190
$map[192] = 195
$map[195] = 192
$map[205] = 206
$map[206] = 205
$map[216] = 217
$map[217] = 221
$map[218] = 216
$map[220] = 218
$map[221] = 220
$map[225] = 226
$map[226] = 225
$map[227] = 231
$map[228] = 230
$map[229] = 228
$map[230] = 229
$map[231] = 227
$map[232] = 233
$map[233] = 234
$map[234] = 232
$map[244] = 245
$map[245] = 246
$map[246] = 244
$map[248] = 251
$map[249] = 248
$map[251] = 249
$map[258] = 259
$map[259] = 258
$map[262] = 263
$map[263] = 264
$map[264] = 265
$map[265] = 262
$map[267] = 269
$map[269] = 267
$map[270] = 271
$map[271] = 270
$map[273] = 274
$map[274] = 273
$map[278] = 280
$map[280] = 278
$map[283] = 288
$map[284] = 285
$map[285] = 287
$map[287] = 283
$map[288] = 284
$map[289] = 292
$map[292] = 289
$map[293] = 295
$map[295] = 293
$map[297] = 298
$map[298] = 297
$map[303] = 304
$map[304] = 303
$map[305] = 307
$map[307] = 308
$map[308] = 305
$map[309] = 311
$map[310] = 309
$map[311] = 310
$map[313] = 314
$map[314] = 317
$map[316] = 313
$map[317] = 316
$map[319] = 320
$map[320] = 319
$map[321] = 322
$map[322] = 321
$map[327] = 328
$map[328] = 327
$map[336] = 337
$map[337] = 336
$map[338] = 340
$map[340] = 341
$map[341] = 338
$map[343] = 345
$map[345] = 343
$map[346] = 347
$map[347] = 346
$map[353] = 354
$map[354] = 355
$map[355] = 353
$map[356] = 358
$map[357] = 356
$map[358] = 359
$map[359] = 357
$map[364] = 366
$map[365] = 364
$map[366] = 365
$map[368] = 369
$map[369] = 368
$map[376] = 377
$map[377] = 376
$map[380] = 381
$map[381] = 380
$map[387] = 389
$map[389] = 387
$map[393] = 397
$map[394] = 393
$map[397] = 394
$map[404] = 405
$map[405] = 404
$map[407] = 412
$map[408] = 409
$map[409] = 408
$map[411] = 407
$map[412] = 411
$map[413] = 416
$map[414] = 415
$map[415] = 414
$map[416] = 418
$map[417] = 413
$map[418] = 417
$map[419] = 420
$map[420] = 419

$affectedRows = @(58, 59, 76, 77, 93, 94, 95, 105, 106, 113, 114, 144, 145, 150, 151, 172, 173, 174, 175, 181, 182, 186, 187, 190, 191, 192, 195, 205, 206, 216, 217, 218, 220, 221, 225, 226, 227, 228, 229, 230, 231, 232, 233, 234, 244, 245, 246, 248, 249, 251, 258, 259, 262, 263, 264, 265, 267, 269, 270, 271, 273, 274, 278, 280, 283, 284, 285, 287, 288, 289, 292, 293, 295, 297, 298, 303, 304, 305, 307, 308, 309, 310, 311, 313, 314, 316, 317, 319, 320, 321, 322, 327, 328, 336, 337, 338, 340, 341, 343, 345, 346, 347, 353, 354, 355, 356, 357, 358, 359, 364, 365, 366, 368, 369, 376, 377, 380, 381, 387, 389, 393, 394, 397, 404, 405, 407, 408, 409, 411, 412, 413, 414, 415, 416, 417, 418, 419, 420)

# Columns: 1=A (ASV_ID) .. 55=BC. Column 6 (F, ASV_rank) is positional
# and must NOT be copied - it stays tied to the row, not the species.
$maxCol = 55
$skipCol = 6

# 1) Snapshot the current (pre-edit) content of every affected row, since the
#    permutation forms cycles and a row can be both a source and a target.
$snapshot = @{}
foreach ($r in $affectedRows) {
    $rowVals = @()
    for ($c = 1; $c -le $maxCol; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each target row full content (A-E, G-BC) from its mapped source row.
foreach ($r in $affectedRows) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    for ($c = 1; $c -le $maxCol; $c++) {
        if ($c -eq $skipCol) { continue }
        $ws.Cells.Item($r, $c).Value = $srcVals[$c - 1]
    }
}

Write-Output "done"